# DISCOVERYACCESS-8155: Add Olin X03 locations to facet mapping.
# Insert three new rows above the existing "OKU Processing" row (row 89)
# for the new Olin Library Room 303 / 403 / 603 facet mappings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A89:A91").EntireRow.Insert()

$ws.Range("A91").Value = "Olin Library Room 603"
$ws.Range("E91").Value = "Olin Library > Room 603"

$ws.Range("A89").Value = "Olin Library Room 303"
$ws.Range("A90").Value = "Olin Library Room 403"

$ws.Range("E90").Value = "Olin Library > Room 403"
$ws.Range("E89").Value = "Olin Library > Room 303"

$null = $ws.Rows(89).EntireRow.Select()
